# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# The worksheet lists overdue-payment periods ("Periodo Mora") for two
# workers. Previously the rows were grouped by worker (all of worker 1's
# periods, then all of worker 2's periods, each block ordered from the
# newest period down to the oldest). The database was refreshed so the
# rows are now interleaved by worker and sorted chronologically by period
# (oldest -> newest): 2002, 2003, ... 2012, 2101, 2102, with each period
# followed immediately by both workers' records.
#
# Column layout (row 15 is the header row):
#   B = Tipo Doc Trabajador (unchanged, always "CC")
#   C = N Doc Trabajador
#   D = Nombre Trabajador
#   E = Periodo Mora
#   F = Valor Mora
#   G = Salario Basico (unchanged, always 877803)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$worker1Doc  = "73290076"
$worker1Name = "PEDRO RAFAEL CARO DE LA HOZ"
$worker2Doc  = "5725752"
$worker2Name = "MIGUEL ANGEL ESPAÑA SAUMETH"

$periods = @("2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012","2101","2102")

$row = 16
foreach ($periodo in $periods) {
    # "Valor Mora" is 25749 only for period 2102, and 35112 for every other period.
    if ($periodo -eq "2102") {
        $valorMora = 25749
    } else {
        $valorMora = 35112
    }

    $ws.Range("C$row").Value = $worker1Doc
    $ws.Range("D$row").Value = $worker1Name
    $ws.Range("E$row").Value = $periodo
    $ws.Range("F$row").Value = $valorMora
    $row = $row + 1

    $ws.Range("C$row").Value = $worker2Doc
    $ws.Range("D$row").Value = $worker2Name
    $ws.Range("E$row").Value = $periodo
    $ws.Range("F$row").Value = $valorMora
    $row = $row + 1
}
